# Apply the MetaData.xlsx edits described by the commit diff.
#
# Summary of changes:
#  - HouseTables:      rows 2-7 gain a "COL4" value in column G;
#                       selection moves to J8 (and the stale topLeftCell
#                       scroll-pin is cleared as a side effect of reselecting).
#  - EducationTables:   rows 2-8 gain StartCode/EndCode/COL1/(COL4|COL5)
#                       data in C/D/F/G, and column B is normalised from a
#                       mix of P3S03/P3S04/P3S05/P3S06/P3S07/P3S08/P3S13
#                       values to the single "P3S09" label for rows 2-14;
#                       selection moves to H6. Once nothing references the
#                       old "P3S07"/"P3S08" shared strings any more they are
#                       dropped automatically on save.
#  - PubWageTable:      selection moves to J21 (no data change).
#  - FoodTables:        rows 2-7 gain COL3/COL4_5/COL6 data in G/J/L;
#                       selection moves to H6.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# HouseTables
# ---------------------------------------------------------------------
$house = $wb.Worksheets.Item("HouseTables")
foreach ($r in 2..7) {
    $house.Cells.Item($r, 7).Value = "COL4"   # column G
}

# ---------------------------------------------------------------------
# EducationTables
# ---------------------------------------------------------------------
$edu = $wb.Worksheets.Item("EducationTables")

foreach ($r in 2..8) {
    $edu.Cells.Item($r, 2).Value = "P3S09"    # column B
    $edu.Cells.Item($r, 3).Value = 72114      # column C
    $edu.Cells.Item($r, 4).Value = 72170      # column D
    $edu.Cells.Item($r, 6).Value = "COL1"     # column F
    if ($r -eq 8) {
        $edu.Cells.Item($r, 7).Value = "COL5" # column G (row 8 only)
    } else {
        $edu.Cells.Item($r, 7).Value = "COL4" # column G
    }
}

foreach ($r in 9..14) {
    $edu.Cells.Item($r, 2).Value = "P3S09"    # column B
}

# ---------------------------------------------------------------------
# FoodTables
# ---------------------------------------------------------------------
$food = $wb.Worksheets.Item("FoodTables")
foreach ($r in 2..7) {
    $food.Cells.Item($r, 7).Value  = "COL3"   # column G
    $food.Cells.Item($r, 10).Value = "COL4_5" # column J
    $food.Cells.Item($r, 12).Value = "COL6"   # column L
}

# ---------------------------------------------------------------------
# Selection / active-cell updates (also clears stale topLeftCell pins)
# ---------------------------------------------------------------------
$house.Range("J8").Select() | Out-Null

$pubWage = $wb.Worksheets.Item("PubWageTable")
$pubWage.Range("J21").Select() | Out-Null

$food.Range("H6").Select() | Out-Null

# EducationTables is reselected last so it remains the active tab,
# matching the workbook's original tabSelected/activeTab state.
$edu.Range("H6").Select() | Out-Null
